# Auto-generated edit script applying cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.243.18"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "1.690.78"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'216.63"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").Value = "'0.523"
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'23.13"
$ws.Range("E8").Value = "  +13.65%  "
$ws.Range("D9").Value = "'0.263"
$ws.Range("E9").Value = "  +4.80%  "
$ws.Range("D10").Value = "'0.0629"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("D11").Value = "'0.0891"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").Value = "1.926.95"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.20"
$ws.Range("E13").Value = "  +2.77%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.678.03"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "'0.555"
$ws.Range("E15").Value = "  +4.99%  "
$ws.Range("D16").Value = "'67.48"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("D17").Value = "27.224.17"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "'238.20"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "0.0₃0747"
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  +3.32%  "
$ws.Range("D23").Value = "'9.69"
$ws.Range("E23").Value = "  +5.40%  "
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("D25").Value = "'148.43"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("D26").Value = "'7.33"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").Value = "'16.56"
$ws.Range("E27").Value = "  +3.06%  "
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  +1.17%  "
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").Value = "1.583.03"
$ws.Range("E32").Value = "  +6.97%  "
$ws.Range("E33").Value = "  +2.71%  "
$ws.Range("E34").Value = "  +2.58%  "
$ws.Range("E35").Value = "  +0.71%  "
$ws.Range("D36").Value = "'0.958"
$ws.Range("E36").Value = "  +6.11%  "
$ws.Range("D37").Value = "'0.605"
$ws.Range("E37").Value = "  +4.00%  "
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "'1.06"
$ws.Range("E40").Value = "  +4.22%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'69.68"
$ws.Range("E41").Value = "  +3.34%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'5.79"
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("D45").Value = "1.836.47"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("D46").Value = "'0.788"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("D47").Value = "'91.40"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").Value = "'1.62"
$ws.Range("E48").Value = "  +6.22%  "
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("D50").Value = "'8.26"
$ws.Range("E50").Value = "  +6.70%  "
$ws.Range("E51").Value = "  +3.49%  "
